$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column stays text (values like "60.758.48" are not real numbers,
# and some like "570.47" would otherwise be auto-converted to numeric by Excel).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '60.760.85'
$ws.Range('E2').Value = '  -2.33%  '
$ws.Range('D3').Value = '2.400.07'
$ws.Range('E3').Value = '  -2.10%  '
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').Value = '570.22'
$ws.Range('E5').Value = '  -1.75%  '
$ws.Range('D6').Value = '139.44'
$ws.Range('E6').Value = '  -2.87%  '
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('D9').Value = '2.379.03'
$ws.Range('E9').Value = '  -2.84%  '
$ws.Range('E10').Value = '  -0.14%  '
$ws.Range('E11').Value = '  +0.24%  '
$ws.Range('E12').Value = '  -2.56%  '
$ws.Range('E13').Value = '  -1.73%  '
$ws.Range('D14').Value = '25.91'
$ws.Range('E14').Value = '  -2.32%  '
$ws.Range('D15').Value = '0.0000169'
$ws.Range('E15').Value = '  -1.93%  '
$ws.Range('D16').Value = '2.800.45'
$ws.Range('D17').Value = '60.721.04'
$ws.Range('E17').Value = '  -2.15%  '
$ws.Range('D18').Value = '2.373.99'
$ws.Range('E18').Value = '  -2.36%  '
$ws.Range('E19').Value = '  -3.31%  '
$ws.Range('D20').Value = '7.20'
$ws.Range('E20').Value = '  +0.68%  '
$ws.Range('D21').Value = '321.43'
$ws.Range('E21').Value = '  -2.41%  '
$ws.Range('D22').Value = '4.01'
$ws.Range('E22').Value = '  -2.01%  '
$ws.Range('D23').Value = '6.07'
$ws.Range('E23').Value = '  +1.26%  '
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('E25').Value = '  -6.85%  '
$ws.Range('D26').Value = '64.28'
$ws.Range('D27').Value = '8.56'
$ws.Range('E27').Value = '  -8.56%  '
$ws.Range('D28').Value = '570.31'
$ws.Range('E28').Value = '  -7.89%  '
$ws.Range('D29').Value = '2.504.24'
$ws.Range('E29').Value = '  -1.46%  '
$ws.Range('D30').Value = '0.0₃0909'
$ws.Range('E30').Value = '  -5.11%  '
$ws.Range('D31').Value = '7.82'
$ws.Range('E31').Value = '  -2.41%  '
$ws.Range('E32').Value = '  -6.24%  '
$ws.Range('E33').Value = '  -3.11%  '
$ws.Range('E34').Value = '  -7.05%  '
$ws.Range('E35').Value = '  +0.31%  '
$ws.Range('E36').Value = '  -5.99%  '
$ws.Range('E37').Value = '  -2.50%  '
$ws.Range('E38').Value = '  -3.78%  '
$ws.Range('E39').Value = '  -1.03%  '
$ws.Range('D40').Value = '146.73'
$ws.Range('E40').Value = '  -2.35%  '
$ws.Range('D41').Value = '5.06'
$ws.Range('E41').Value = '  -4.80%  '
$ws.Range('E42').Value = '  +0.14%  '
$ws.Range('D43').Value = '41.54'
$ws.Range('E43').Value = '  -2.28%  '
$ws.Range('E44').Value = '  -4.99%  '
$ws.Range('E45').Value = '  -5.72%  '
$ws.Range('D46').Value = '0.0₆0283'
$ws.Range('E46').Value = '  +17.09%  '
$ws.Range('D47').Value = '140.07'
$ws.Range('E47').Value = '  -2.24%  '
$ws.Range('E48').Value = '  -4.04%  '
$ws.Range('D49').Value = '0.583'
$ws.Range('E49').Value = '  -3.60%  '
$ws.Range('E50').Value = '  -4.16%  '
$ws.Range('D51').Value = '19.25'
$ws.Range('E51').Value = '  -1.59%  '
